# Search test cases modified
# - Set Runmode (column D) to "Y" for all test case rows (2-20) on "Test Cases" sheet
# - Append "|OPQA-511" to the JIRA ID (column B) of the row whose JIRA ID is "OPQA-496"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 4).Value2 = "Y"

    $jira = $ws.Cells.Item($r, 2).Value2
    if ($jira -eq "OPQA-496") {
        $ws.Cells.Item($r, 2).Value2 = "OPQA-496|OPQA-511"
    }
}
